$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# --- Upload the traditional matrix-multiplication timings ---

# Int results (rows 3-12, column B = "Traditional")
$ws.Range("B3").Value  = 0.58199999999999996
$ws.Range("B4").Value  = 4.9720000000000004
$ws.Range("B5").Value  = 16.312999999999999
$ws.Range("B6").Value  = 37.896999999999998
$ws.Range("B7").Value  = 73.917000000000002
$ws.Range("B8").Value  = 127.797
$ws.Range("B9").Value  = 202.16900000000001
$ws.Range("B10").Value = 303.73099999999999
$ws.Range("B11").Value = 436.89699999999999
$ws.Range("B12").Value = 595.43100000000004

# Float results (rows 24-33, column B = "Traditional")
$ws.Range("B24").Value = 0.61
$ws.Range("B25").Value = 4.8360000000000003
$ws.Range("B26").Value = 16.401
$ws.Range("B27").Value = 38.927999999999997
$ws.Range("B28").Value = 76.033000000000001
$ws.Range("B29").Value = 131.85900000000001
$ws.Range("B30").Value = 208.809
$ws.Range("B31").Value = 312.774
$ws.Range("B32").Value = 443.23700000000002
$ws.Range("B33").Value = 609.21199999999999

# --- Add the "Int Results" chart (mirrors the existing float chart) ---

$shp = $ws.Shapes.AddChart2(-1, 51)
$chart = $shp.Chart
$chart.SetSourceData($ws.Range("A2:C12"))
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Int Results"
$chart.HasLegend = $true
$chart.Legend.Position = -4107

# Position/size the new chart to sit above the existing "Float Results"
# chart, spanning columns E:M, rows 1:15 (matches the author's layout)
$shp.Left = 490.431640625
$shp.Top = 0
$shp.Width = 443.5
$shp.Height = 216

# Leave the active selection where the author's session ended up
$ws.Range("Q20").Select()
